$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("U28:U120").Select()
